$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H132").Value = 3970065.2
$ws.Range("I132").Value = 4610101.5
$ws.Range("J132").Value = 1839.8
$ws.Range("K132").Value = 13830304.5
$ws.Range("L132").Value = 5519.4
$ws.Range("M132").Value = -13827774.5
$ws.Range("N132").Value = -10579.4
$ws.Range("H137").Value = 1776.5454
$ws.Range("I137").Value = 1563.7646
$ws.Range("K137").Value = 4691.293799999999
$ws.Range("M137").Value = -2141.293799999999
$ws.Range("H141").Value = 2734.3914
$ws.Range("I141").Value = 1535.909
$ws.Range("J141").Value = 3833
$ws.Range("K141").Value = 4607.727000000001
$ws.Range("L141").Value = 11499
$ws.Range("M141").Value = 572.2729999999992
$ws.Range("N141").Value = -21859
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2856.889
$ws.Range("I122").Value = 1853
$ws.Range("K122").Value = 5559
$ws.Range("M122").Value = -3109
$ws.Range("H131").Value = 41000
$ws.Range("J131").Value = 41000
$ws.Range("L131").Value = 41000
$ws.Range("N131").Value = -51080
$ws.Range("H132").Value = 4102.3057
$ws.Range("I132").Value = 5464.609
$ws.Range("J132").Value = 1692.0769
$ws.Range("K132").Value = 16393.827
$ws.Range("L132").Value = 5076.2307
$ws.Range("M132").Value = -13863.827
$ws.Range("N132").Value = -10136.2307
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6668123.5
$ws.Range("I31").Value = 1652.6666
$ws.Range("J31").Value = 22223222
$ws.Range("K31").Value = 1652.6666
$ws.Range("L31").Value = 22223222
$ws.Range("M31").Value = -1357.6666
$ws.Range("N31").Value = -22223812
$ws.Range("H34").Value = 6668123.5
$ws.Range("I34").Value = 1652.6666
$ws.Range("J34").Value = 22223222
$ws.Range("K34").Value = 1652.6666
$ws.Range("L34").Value = 22223222
$ws.Range("M34").Value = -1450.6666
$ws.Range("N34").Value = -22223626
$ws.Range("H64").Value = 50000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 50000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H107").Value = 577.3939
$ws.Range("I107").Value = 538.8946999999999
$ws.Range("J107").Value = 629.6429000000001
$ws.Range("K107").Value = 538.8946999999999
$ws.Range("L107").Value = 629.6429000000001
$ws.Range("M107").Value = 1381.1053
$ws.Range("N107").Value = -4469.6429
$ws.Range("H141").Value = 55000
$ws.Range("I141").Value = 50000
$ws.Range("J141").Value = 60000
$ws.Range("K141").Value = 50000
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = -44820
$ws.Range("N141").Value = -70360
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 9000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -10872
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 27000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -36360
$ws.Range("H131").Value = 4449604.5
$ws.Range("J131").Value = 5291817
$ws.Range("L131").Value = 15875451
$ws.Range("N131").Value = -15885531
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1195.0416
$ws.Range("I102").Value = 1183.2106
$ws.Range("J102").Value = 1240
$ws.Range("K102").Value = 1183.2106
$ws.Range("L102").Value = 1240
$ws.Range("M102").Value = 438.7893999999999
$ws.Range("N102").Value = -4484
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 251751
$ws.Range("I40").Value = 334668
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 334668
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -334532
$ws.Range("N40").Value = -3272
$ws.Range("H122").Value = 3100.4707
$ws.Range("I122").Value = 3213.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 9640.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -7190.5
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 7710.067
$ws.Range("I132").Value = 18609.9
$ws.Range("J132").Value = 2260.15
$ws.Range("K132").Value = 55829.7
$ws.Range("L132").Value = 6780.450000000001
$ws.Range("M132").Value = -53299.7
$ws.Range("N132").Value = -11840.45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2113
$ws.Range("I81").Value = 2500
$ws.Range("J81").Value = 2027
$ws.Range("K81").Value = 5000
$ws.Range("L81").Value = 4054
$ws.Range("M81").Value = -3939
$ws.Range("N81").Value = -6176
$ws.Range("H84").Value = 2113
$ws.Range("I84").Value = 2500
$ws.Range("J84").Value = 2027
$ws.Range("K84").Value = 25000
$ws.Range("L84").Value = 20270
$ws.Range("M84").Value = -19696
$ws.Range("N84").Value = -30878
$ws.Range("H126").Value = 50652
$ws.Range("I126").Value = 50652
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 151956
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -149486
$ws.Range("N126").ClearContents()
$ws.Range("H130").Value = 49996.668
$ws.Range("J130").Value = 49996.668
$ws.Range("L130").Value = 49996.668
$ws.Range("N130").Value = -60036.668
$ws.Range("H132").Value = 1001.3889
$ws.Range("I132").Value = 867.375
$ws.Range("J132").Value = 2073.5
$ws.Range("K132").Value = 2602.125
$ws.Range("L132").Value = 6220.5
$ws.Range("M132").Value = -72.125
$ws.Range("N132").Value = -11280.5
$ws.Range("H136").Value = 6577.5454
$ws.Range("I136").Value = 8218
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 24654
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -22104
$ws.Range("N136").Value = -8100
